# Overworld stuff, more UI art: add climate types and region entries to the
# Language sheet (Climate block gets 6 new rows, Region block gets 9 new rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert 6 new rows for additional climate types, right before the
#     existing "region_title" row (currently row 34). ---
$ws.Range("A34:A39").EntireRow.Insert()

$ws.Cells.Item(34,1).Value = "climate_tropical"
$ws.Cells.Item(34,2).Value = "Tropical"

$ws.Cells.Item(35,1).Value = "climate_oceanic"
$ws.Cells.Item(35,2).Value = "Oceanic"

$ws.Cells.Item(36,1).Value = "climate_desert"
$ws.Cells.Item(36,2).Value = "Desert"

$ws.Cells.Item(37,1).Value = "climate_tundra"
$ws.Cells.Item(37,2).Value = "Tundra"

$ws.Cells.Item(38,1).Value = "climate_highland"
$ws.Cells.Item(38,2).Value = "Highland"

$ws.Cells.Item(39,1).Value = "climate_mediterranean"
$ws.Cells.Item(39,2).Value = "Mediterranean"

# After the insert above, "region_title" is now row 40 and "region_NA" is
# row 41. Insert 9 new rows right after "region_NA" for the new regions.
$ws.Range("A42:A50").EntireRow.Insert()

$ws.Cells.Item(42,1).Value = "region_PH"
$ws.Cells.Item(42,2).Value = "Northern Luzon, Philippines"

$ws.Cells.Item(43,1).Value = "region_GB"
$ws.Cells.Item(43,2).Value = "British Isles"

$ws.Cells.Item(44,1).Value = "region_MG"
$ws.Cells.Item(44,2).Value = "Madagascar"

$ws.Cells.Item(45,1).Value = "region_EG"
$ws.Cells.Item(45,2).Value = "Egypt"

$ws.Cells.Item(46,1).Value = "region_GL"
$ws.Cells.Item(46,2).Value = "Greenland"

$ws.Cells.Item(47,1).Value = "region_BR"
$ws.Cells.Item(47,2).Value = "Brazil"

$ws.Cells.Item(48,1).Value = "region_CL"
$ws.Cells.Item(48,2).Value = "Andes Mountains, Chile"

$ws.Cells.Item(49,1).Value = "region_IT"
$ws.Cells.Item(49,2).Value = "Italy"

$ws.Cells.Item(50,1).Value = "region_AU"
$ws.Cells.Item(50,2).Value = "Australian Outback"

# Match the final selection state from the authored edit.
$ws.Range("B50").Select()
